$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I7").Value = "sd"
$ws.Range("J7").Value = "Statement-non-opinion"
$ws.Range("I8").Value = "sd"
$ws.Range("J8").Value = "Statement-non-opinion"
$ws.Range("I16").Value = "sd"
$ws.Range("J16").Value = "Statement-non-opinion"
$ws.Range("I40").Value = "sv"
$ws.Range("J40").Value = "Statement-opinion"
$ws.Range("I48").Value = "sd"
$ws.Range("J48").Value = "Statement-non-opinion"
$ws.Range("I49").Value = "aa"
$ws.Range("J49").Value = "Agree/Accept"
$ws.Range("I63").Value = "ba"
$ws.Range("J63").Value = "Appreciation"
$ws.Range("I71").Value = "sd"
$ws.Range("J71").Value = "Statement-non-opinion"
$ws.Range("I79").Value = "sv"
$ws.Range("J79").Value = "Statement-opinion"
$ws.Range("I96").Value = "sd"
$ws.Range("J96").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "sd"
$ws.Range("J117").Value = "Statement-non-opinion"
$ws.Range("I134").Value = "sd"
$ws.Range("J134").Value = "Statement-non-opinion"
$ws.Range("I141").Value = "sv"
$ws.Range("J141").Value = "Statement-opinion"
$ws.Range("I150").Value = "sd"
$ws.Range("J150").Value = "Statement-non-opinion"
$ws.Range("I154").Value = "qy"
$ws.Range("J154").Value = "Yes-No-Question"
$ws.Range("I158").Value = "sv"
$ws.Range("J158").Value = "Statement-opinion"
$ws.Range("I173").Value = "sd"
$ws.Range("J173").Value = "Statement-non-opinion"
$ws.Range("I174").Value = "sd"
$ws.Range("J174").Value = "Statement-non-opinion"
$ws.Range("I193").Value = "aa"
$ws.Range("J193").Value = "Agree/Accept"
$ws.Range("I211").Value = "sd"
$ws.Range("J211").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "sv"
$ws.Range("J214").Value = "Statement-opinion"
$ws.Range("I222").Value = "sd"
$ws.Range("J222").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "ba"
$ws.Range("J243").Value = "Appreciation"
$ws.Range("I258").Value = "aa"
$ws.Range("J258").Value = "Agree/Accept"
$ws.Range("I269").Value = "sv"
$ws.Range("J269").Value = "Statement-opinion"
$ws.Range("I271").Value = "sv"
$ws.Range("J271").Value = "Statement-opinion"
$ws.Range("I278").Value = "sv"
$ws.Range("J278").Value = "Statement-opinion"
$ws.Range("I282").Value = "ba"
$ws.Range("J282").Value = "Appreciation"
$ws.Range("I303").Value = "sv"
$ws.Range("J303").Value = "Statement-opinion"
